$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Populate Sheet1 with the new login-table test data first, so the ---
# --- shared-string table is built in the same order as the source file ---
# (row 2 first, then the header row, then row 3).
$ws1.Range("A2").Value = "Disney"
$ws1.Range("B2").Value = "WaltWhite"
$ws1.Range("C2").Value = "password1"

$ws1.Range("A1").Value = "CompanyId"
$ws1.Range("B1").Value = "inputID"
$ws1.Range("C1").Value = "Password"

$ws1.Range("A3").Value = "sree"
$ws1.Range("B3").Value = "sree1"
$ws1.Range("C3").Value = "password1"

# --- Add Sheet2 as a copy of Sheet1 before the header highlight / column ---
# --- widths are applied, then shrink it back down to a single column ---
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Sheet2"
$ws2.Range("A1").Value = "Browser"
$ws2.Range("A2").Value = "chrome"
$ws2.Range("B1:C2").Value = ""
$ws2.Range("A3:C3").Value = ""

# --- Finish formatting Sheet1: highlight header row, size columns ---
$ws1.Range("A1:C1").Interior.ColorIndex = 6
$ws1.Columns("A").ColumnWidth = 10.1
$ws1.Columns("B").ColumnWidth = 9.65
$ws1.Columns("C").ColumnWidth = 9.65

$ws1.Range("I16").Select() | Out-Null
$ws2.Range("B4").Select() | Out-Null
